# Auto-generated script to update market-price / profit columns (H:N)
# across the Kujata_Profits leve sheets, per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value2 = 325.5
$ws.Range("J12").Value2 = 550.5
$ws.Range("L12").Value2 = 550.5
$ws.Range("N12").Value2 = -890.5
$ws.Range("H64").Value2 = 3933.5557
$ws.Range("I64").Value2 = 3968.4
$ws.Range("J64").Value2 = 3890
$ws.Range("K64").Value2 = 3968.4
$ws.Range("L64").Value2 = 3890
$ws.Range("M64").Value2 = -3720.4
$ws.Range("N64").Value2 = -4386
$ws.Range("H67").Value2 = 3933.5557
$ws.Range("I67").Value2 = 3968.4
$ws.Range("J67").Value2 = 3890
$ws.Range("K67").Value2 = 3968.4
$ws.Range("L67").Value2 = 3890
$ws.Range("M67").Value2 = -3110.4
$ws.Range("N67").Value2 = -5606
$ws.Range("H74").Value2 = 4751.222
$ws.Range("I74").Value2 = 2800
$ws.Range("K74").Value2 = 2800
$ws.Range("M74").Value2 = -1864
$ws.Range("H75").Value2 = 20431.4
$ws.Range("J75").Value2 = 20431.4
$ws.Range("L75").Value2 = 20431.4
$ws.Range("N75").Value2 = -22303.4
$ws.Range("H76").Value2 = 5250
$ws.Range("I76").Value2 = 3000
$ws.Range("J76").Value2 = 7500
$ws.Range("K76").Value2 = 3000
$ws.Range("L76").Value2 = 7500
$ws.Range("M76").Value2 = -2685
$ws.Range("N76").Value2 = -8130
$ws.Range("H77").Value2 = 4751.222
$ws.Range("I77").Value2 = 2800
$ws.Range("K77").Value2 = 14000
$ws.Range("M77").Value2 = -9320
$ws.Range("H78").Value2 = 20431.4
$ws.Range("J78").Value2 = 20431.4
$ws.Range("L78").Value2 = 61294.2
$ws.Range("N78").Value2 = -70654.20000000001
$ws.Range("H79").Value2 = 5250
$ws.Range("I79").Value2 = 3000
$ws.Range("J79").Value2 = 7500
$ws.Range("K79").Value2 = 3000
$ws.Range("L79").Value2 = 7500
$ws.Range("M79").Value2 = -1908
$ws.Range("N79").Value2 = -9684
$ws.Range("H96").Value2 = 2144.5557
$ws.Range("I96").Value2 = 3178.4
$ws.Range("K96").Value2 = 9535.200000000001
$ws.Range("M96").Value2 = -8162.200000000001
$ws.Range("H100").Value2 = 2134.5454
$ws.Range("J100").Value2 = 2525.7144
$ws.Range("L100").Value2 = 2525.7144
$ws.Range("N100").Value2 = -3607.7144
$ws.Range("H106").Value2 = 7504.1577
$ws.Range("I106").Value2 = 7851.6113
$ws.Range("J106").Value2 = 1250
$ws.Range("K106").Value2 = 7851.6113
$ws.Range("L106").Value2 = 1250
$ws.Range("M106").Value2 = -7220.6113
$ws.Range("N106").Value2 = -2512
$ws.Range("H116").Value2 = 3408.2222
$ws.Range("I116").Value2 = 3524.8572
$ws.Range("J116").Value2 = 3000
$ws.Range("K116").Value2 = 3524.8572
$ws.Range("L116").Value2 = 3000
$ws.Range("M116").Value2 = -82.85719999999992
$ws.Range("N116").Value2 = -9884

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 10148.43
$ws.Range("I32").Value2 = 7522.9727
$ws.Range("J32").Value2 = 19731.35
$ws.Range("K32").Value2 = 7522.9727
$ws.Range("L32").Value2 = 19731.35
$ws.Range("M32").Value2 = -7235.9727
$ws.Range("N32").Value2 = -20305.35
$ws.Range("H63").Value2 = 32260348
$ws.Range("J63").Value2 = 200002200
$ws.Range("L63").Value2 = 200002200
$ws.Range("N63").Value2 = -200003572
$ws.Range("H66").Value2 = 32260348
$ws.Range("J66").Value2 = 200002200
$ws.Range("L66").Value2 = 1000011000
$ws.Range("N66").Value2 = -1000017864
$ws.Range("H97").Value2 = 6643
$ws.Range("I97").Value2 = 826.6667
$ws.Range("J97").Value2 = 20602.2
$ws.Range("K97").Value2 = 826.6667
$ws.Range("L97").Value2 = 20602.2
$ws.Range("M97").Value2 = -330.6667
$ws.Range("N97").Value2 = -21594.2
$ws.Range("H102").Value2 = 8773465
$ws.Range("I102").Value2 = 10417802
$ws.Range("K102").Value2 = 10417802
$ws.Range("M102").Value2 = -10416180
$ws.Range("H117").Value2 = 63333.332
$ws.Range("J117").Value2 = 63333.332
$ws.Range("L117").Value2 = 63333.332
$ws.Range("N117").Value2 = -72511.33199999999
$ws.Range("H132").Value2 = 4545.136
$ws.Range("I132").Value2 = 4721.615
$ws.Range("J132").Value2 = 4290.222
$ws.Range("K132").Value2 = 14164.845
$ws.Range("L132").Value2 = 12870.666
$ws.Range("M132").Value2 = -11634.845
$ws.Range("N132").Value2 = -17930.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 3845.3635
$ws.Range("I20").Value2 = 3462.5
$ws.Range("K20").Value2 = 3462.5
$ws.Range("M20").Value2 = -3215.5
$ws.Range("H80").Value2 = 568.7059
$ws.Range("I80").Value2 = 218.42857
$ws.Range("K80").Value2 = 218.42857
$ws.Range("M80").Value2 = 779.57143
$ws.Range("H83").Value2 = 568.7059
$ws.Range("I83").Value2 = 218.42857
$ws.Range("K83").Value2 = 1092.14285
$ws.Range("M83").Value2 = 3899.85715
$ws.Range("H86").Value2 = 3792.0435
$ws.Range("I86").Value2 = 3882.5881
$ws.Range("K86").Value2 = 3882.5881
$ws.Range("M86").Value2 = -2759.5881
$ws.Range("H89").Value2 = 3792.0435
$ws.Range("I89").Value2 = 3882.5881
$ws.Range("K89").Value2 = 19412.9405
$ws.Range("M89").Value2 = -13796.9405
$ws.Range("H96").Value2 = 31000
$ws.Range("I96").Value2 = 0
$ws.Range("J96").Value2 = 31000
$ws.Range("K96").Value2 = 0
$ws.Range("L96").Value2 = 31000
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value2 = -36492
$ws.Range("H99").Value2 = 90910300
$ws.Range("I99").Value2 = 125001050
$ws.Range("K99").Value2 = 125001050
$ws.Range("M99").Value2 = -124999552
$ws.Range("H105").Value2 = 144273620
$ws.Range("I105").Value2 = 168318880
$ws.Range("K105").Value2 = 168318880
$ws.Range("M105").Value2 = -168317133
$ws.Range("H134").Value2 = 4058.5151
$ws.Range("I134").Value2 = 961.6923
$ws.Range("J134").Value2 = 15561
$ws.Range("K134").Value2 = 2885.0769
$ws.Range("L134").Value2 = 46683
$ws.Range("M134").Value2 = -350.0769
$ws.Range("N134").Value2 = -51753

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 1565.6038
$ws.Range("I31").Value2 = 1383.0465
$ws.Range("J31").Value2 = 2350.6
$ws.Range("K31").Value2 = 1383.0465
$ws.Range("L31").Value2 = 2350.6
$ws.Range("M31").Value2 = -1088.0465
$ws.Range("N31").Value2 = -2940.6
$ws.Range("H34").Value2 = 1565.6038
$ws.Range("I34").Value2 = 1383.0465
$ws.Range("J34").Value2 = 2350.6
$ws.Range("K34").Value2 = 1383.0465
$ws.Range("L34").Value2 = 2350.6
$ws.Range("M34").Value2 = -1181.0465
$ws.Range("N34").Value2 = -2754.6
$ws.Range("H99").Value2 = 2025895.2
$ws.Range("I99").Value2 = 3760891.5
$ws.Range("J99").Value2 = 1733.1666
$ws.Range("K99").Value2 = 3760891.5
$ws.Range("L99").Value2 = 1733.1666
$ws.Range("M99").Value2 = -3759393.5
$ws.Range("N99").Value2 = -4729.1666
$ws.Range("H126").Value2 = 2025895.2
$ws.Range("I126").Value2 = 3760891.5
$ws.Range("J126").Value2 = 1733.1666
$ws.Range("K126").Value2 = 11282674.5
$ws.Range("L126").Value2 = 5199.4998
$ws.Range("M126").Value2 = -11280204.5
$ws.Range("N126").Value2 = -10139.4998
$ws.Range("H141").Value2 = 284664.2
$ws.Range("J141").Value2 = 284664.2
$ws.Range("L141").Value2 = 284664.2
$ws.Range("N141").Value2 = -295024.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value2 = 1099
$ws.Range("J97").Value2 = 1099
$ws.Range("L97").Value2 = 3297
$ws.Range("N97").Value2 = -4289
$ws.Range("H98").Value2 = 479.33334
$ws.Range("I98").Value2 = 339.25
$ws.Range("J98").Value2 = 591.4
$ws.Range("K98").Value2 = 1017.75
$ws.Range("L98").Value2 = 1774.2
$ws.Range("M98").Value2 = 480.25
$ws.Range("N98").Value2 = -4770.2
$ws.Range("H107").Value2 = 10637.4
$ws.Range("I107").Value2 = 595
$ws.Range("K107").Value2 = 1785
$ws.Range("M107").Value2 = 135
$ws.Range("H113").Value2 = 692.2917
$ws.Range("J113").Value2 = 748.14813
$ws.Range("L113").Value2 = 2244.44439
$ws.Range("N113").Value2 = -6584.444390000001
$ws.Range("H131").Value2 = 31297806
$ws.Range("I131").Value2 = 76923460
$ws.Range("J131").Value2 = 80256
$ws.Range("K131").Value2 = 230770380
$ws.Range("L131").Value2 = 240768
$ws.Range("M131").Value2 = -230765340
$ws.Range("N131").Value2 = -250848

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 64288860
$ws.Range("I70").Value2 = 41670336
$ws.Range("K70").Value2 = 41670336
$ws.Range("M70").Value2 = -41670066
$ws.Range("H73").Value2 = 64288860
$ws.Range("I73").Value2 = 41670336
$ws.Range("K73").Value2 = 41670336
$ws.Range("M73").Value2 = -41669400
$ws.Range("H80").Value2 = 3744.3333
$ws.Range("I80").Value2 = 1800
$ws.Range("J80").Value2 = 4299.857
$ws.Range("K80").Value2 = 1800
$ws.Range("L80").Value2 = 4299.857
$ws.Range("M80").Value2 = -802
$ws.Range("N80").Value2 = -6295.857
$ws.Range("H83").Value2 = 3744.3333
$ws.Range("I83").Value2 = 1800
$ws.Range("J83").Value2 = 4299.857
$ws.Range("K83").Value2 = 9000
$ws.Range("L83").Value2 = 21499.285
$ws.Range("M83").Value2 = -4008
$ws.Range("N83").Value2 = -31483.285
$ws.Range("H97").Value2 = 1159
$ws.Range("I97").Value2 = 1000.8889
$ws.Range("K97").Value2 = 1000.8889
$ws.Range("M97").Value2 = -504.8889
$ws.Range("H126").Value2 = 2022.3077
$ws.Range("I126").Value2 = 1621.1111
$ws.Range("J126").Value2 = 2925
$ws.Range("K126").Value2 = 4863.3333
$ws.Range("L126").Value2 = 8775
$ws.Range("M126").Value2 = -2393.3333
$ws.Range("N126").Value2 = -13715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 847.6
$ws.Range("I16").Value2 = 847.6
$ws.Range("K16").Value2 = 847.6
$ws.Range("M16").Value2 = -677.6
$ws.Range("H55").Value2 = 352.78262
$ws.Range("I55").Value2 = 263.89474
$ws.Range("J55").Value2 = 775
$ws.Range("K55").Value2 = 263.89474
$ws.Range("L55").Value2 = 775
$ws.Range("M55").Value2 = -90.89474000000001
$ws.Range("N55").Value2 = -1121

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 523.25
$ws.Range("I107").Value2 = 531.3333
$ws.Range("K107").Value2 = 1593.9999
$ws.Range("M107").Value2 = 326.0001
$ws.Range("H121").Value2 = 30000
$ws.Range("J121").Value2 = 30000
$ws.Range("L121").Value2 = 30000
$ws.Range("N121").Value2 = -33494
$ws.Range("H124").Value2 = 30125
$ws.Range("J124").Value2 = 30125
$ws.Range("L124").Value2 = 30125
$ws.Range("N124").Value2 = -39945
